# Auto-generated: apply 2024-07-11 data update across Citywide Totals,
# By Neighborhood summary, and per-neighborhood sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("K2").Value = 4116
$ws.Range("J3").Value = 8079
$ws.Range("K3").Value = 4225
$ws.Range("K5").Value = 304
$ws.Range("K6").Value = 4729
$ws.Range("J7").Value = 29292
$ws.Range("K7").Value = 14221

$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("K7").Value = 415
$ws.Range("K8").Value = 960
$ws.Range("K11").Value = 279
$ws.Range("K14").Value = 76
$ws.Range("K15").Value = 146
$ws.Range("K17").Value = 27
$ws.Range("K18").Value = 99
$ws.Range("K20").Value = 318
$ws.Range("K21").Value = 42
$ws.Range("K23").Value = 145
$ws.Range("K28").Value = 6
$ws.Range("K29").Value = 746
$ws.Range("K31").Value = 154
$ws.Range("K33").Value = 591
$ws.Range("K36").Value = 179
$ws.Range("K37").Value = 483
$ws.Range("K42").Value = 515
$ws.Range("K44").Value = 130
$ws.Range("K47").Value = 83
$ws.Range("K48").Value = 181
$ws.Range("K49").Value = 84
$ws.Range("K50").Value = 76
$ws.Range("K51").Value = 176
$ws.Range("K53").Value = 192
$ws.Range("K55").Value = 161
$ws.Range("J63").Value = 107
$ws.Range("K63").Value = 48
$ws.Range("K65").Value = 327
$ws.Range("K67").Value = 551
$ws.Range("K77").Value = 100
$ws.Range("K79").Value = 367
$ws.Range("K83").Value = 301
$ws.Range("K84").Value = 103
$ws.Range("K85").Value = 636
$ws.Range("K87").Value = 22
$ws.Range("K88").Value = 163
$ws.Range("K89").Value = 201
$ws.Range("K91").Value = 159
$ws.Range("K92").Value = 50
$ws.Range("K94").Value = 178
$ws.Range("K95").Value = 243
$ws.Range("K97").Value = 122
$ws.Range("K99").Value = 243
$ws.Range("J101").Value = 29292
$ws.Range("K101").Value = 14221

$ws = $wb.Worksheets.Item("Bridgeport")
$ws.Range("K3").Value = 16
$ws.Range("K7").Value = 76

$ws = $wb.Worksheets.Item("Auburn Gresham")
$ws.Range("K3").Value = 136
$ws.Range("K6").Value = 101
$ws.Range("K7").Value = 415

$ws = $wb.Worksheets.Item("Belmont Cragin")
$ws.Range("K3").Value = 71
$ws.Range("K7").Value = 279

$ws = $wb.Worksheets.Item("Uptown")
$ws.Range("K2").Value = 54
$ws.Range("K3").Value = 57
$ws.Range("K7").Value = 201

$ws = $wb.Worksheets.Item("South Shore")
$ws.Range("K2").Value = 224
$ws.Range("K7").Value = 636

$ws = $wb.Worksheets.Item("Logan Square")
$ws.Range("K6").Value = 91
$ws.Range("K7").Value = 192

$ws = $wb.Worksheets.Item("Austin")
$ws.Range("K3").Value = 287
$ws.Range("K7").Value = 960

$ws = $wb.Worksheets.Item("South Chicago")
$ws.Range("K6").Value = 70
$ws.Range("K7").Value = 301

$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("K2").Value = 160
$ws.Range("K3").Value = 223
$ws.Range("K6").Value = 170
$ws.Range("K7").Value = 591

$ws = $wb.Worksheets.Item("West Pullman")
$ws.Range("K3").Value = 85
$ws.Range("K7").Value = 243

$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Range("K2").Value = 131
$ws.Range("K3").Value = 160
$ws.Range("K5").Value = 24
$ws.Range("K6").Value = 144
$ws.Range("K7").Value = 483

$ws = $wb.Worksheets.Item("New City")
$ws.Range("K5").Value = 8
$ws.Range("K6").Value = 131
$ws.Range("K7").Value = 327

$ws = $wb.Worksheets.Item("Woodlawn")
$ws.Range("K6").Value = 58
$ws.Range("K7").Value = 243

$ws = $wb.Worksheets.Item("Gage Park")
$ws.Range("K3").Value = 37
$ws.Range("K7").Value = 154

$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Range("K6").Value = 163
$ws.Range("K7").Value = 551

$ws = $wb.Worksheets.Item("South Deering")
$ws.Range("K4").Value = 7
$ws.Range("K6").Value = 23
$ws.Range("K7").Value = 103

$ws = $wb.Worksheets.Item("Lincoln Park")
$ws.Range("K6").Value = 48
$ws.Range("K7").Value = 84

$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("K2").Value = 211
$ws.Range("K3").Value = 268
$ws.Range("K7").Value = 746

$ws = $wb.Worksheets.Item("Lake View")
$ws.Range("K2").Value = 24
$ws.Range("K7").Value = 181

$ws = $wb.Worksheets.Item("Irving Park")
$ws.Range("K4").Value = 7
$ws.Range("K7").Value = 130

$ws = $wb.Worksheets.Item("Humboldt Park")
$ws.Range("K2").Value = 140
$ws.Range("K4").Value = 21
$ws.Range("K6").Value = 186
$ws.Range("K7").Value = 515

$ws = $wb.Worksheets.Item("Lower West Side")
$ws.Range("K4").Value = 7
$ws.Range("K6").Value = 59
$ws.Range("K7").Value = 161

$ws = $wb.Worksheets.Item("Douglas")
$ws.Range("K2").Value = 42
$ws.Range("K3").Value = 54
$ws.Range("K6").Value = 34
$ws.Range("K7").Value = 145

$ws = $wb.Worksheets.Item("Washington Park")
$ws.Range("K2").Value = 40
$ws.Range("K3").Value = 74
$ws.Range("K7").Value = 159

$ws = $wb.Worksheets.Item("Chinatown")
$ws.Range("K6").Value = 23
$ws.Range("K7").Value = 42

$ws = $wb.Worksheets.Item("Roseland")
$ws.Range("K3").Value = 120
$ws.Range("K5").Value = 14
$ws.Range("K7").Value = 367

$ws = $wb.Worksheets.Item("Chicago Lawn")
$ws.Range("K2").Value = 109
$ws.Range("K7").Value = 318

$ws = $wb.Worksheets.Item("Calumet Heights")
$ws.Range("K3").Value = 30
$ws.Range("K6").Value = 25
$ws.Range("K7").Value = 99

$ws = $wb.Worksheets.Item("Burnside")
$ws.Range("K2").Value = 12
$ws.Range("K7").Value = 27

$ws = $wb.Worksheets.Item("Grand Boulevard")
$ws.Range("K5").Value = 2
$ws.Range("K7").Value = 179

$ws = $wb.Worksheets.Item("West Loop")
$ws.Range("K6").Value = 77
$ws.Range("K7").Value = 178

$ws = $wb.Worksheets.Item("Kenwood")
$ws.Range("K4").Value = 5
$ws.Range("K7").Value = 83

$ws = $wb.Worksheets.Item("Brighton Park")
$ws.Range("K3").Value = 39
$ws.Range("K7").Value = 146

$ws = $wb.Worksheets.Item("Lincoln Square")
$ws.Range("K2").Value = 17
$ws.Range("K7").Value = 76

$ws = $wb.Worksheets.Item("West Town")
$ws.Range("K6").Value = 75
$ws.Range("K7").Value = 122

$ws = $wb.Worksheets.Item("West Elsdon")
$ws.Range("K2").Value = 16
$ws.Range("K7").Value = 50

$ws = $wb.Worksheets.Item("United Center")
$ws.Range("K3").Value = 52
$ws.Range("K7").Value = 163

$ws = $wb.Worksheets.Item("Little Italy, UIC")
$ws.Range("K2").Value = 47
$ws.Range("K6").Value = 58
$ws.Range("K7").Value = 176

$ws = $wb.Worksheets.Item("Riverdale")
$ws.Range("K6").Value = 14
$ws.Range("K7").Value = 100

$ws = $wb.Worksheets.Item("Ukrainian Village")
$ws.Range("K6").Value = 12
$ws.Range("K7").Value = 22

$ws = $wb.Worksheets.Item("Edison Park")
$ws.Range("K2").Value = 5
$ws.Range("K7").Value = 6
